$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 974.2
$ws.Range("K5").Value = 725.4
$ws.Range("K6").Value = 1236.3

# Match the formatting of the preceding column (J) for each row
$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial(-4122)

$ws.Range("J5").Copy()
$ws.Range("K5").PasteSpecial(-4122)

$ws.Range("J6").Copy()
$ws.Range("K6").PasteSpecial(-4122)

$excel.CutCopyMode = 0
